# Update "想去人数" (want-to-go count) values in column F across the
# 展览 (Exhibition), 本地生活 (Local life) and 全部类型 (All types) sheets.
# 演出 (Performance) sheet has no changes in this update.

$wb = $excel.ActiveWorkbook

# ---- 展览 sheet ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3762
$ws.Range("F5").Value = 3762
$ws.Range("F6").Value = 295
$ws.Range("F7").Value = 5301
$ws.Range("F8").Value = 590
$ws.Range("F9").Value = 426
$ws.Range("F10").Value = 229
$ws.Range("F11").Value = 1061
$ws.Range("F13").Value = 141
$ws.Range("F14").Value = 50
$ws.Range("F15").Value = 731
$ws.Range("F16").Value = 364
$ws.Range("F17").Value = 47
$ws.Range("F19").Value = 179
$ws.Range("F22").Value = 6067
$ws.Range("F23").Value = 6067
$ws.Range("F25").Value = 44
$ws.Range("F27").Value = 6826
$ws.Range("F30").Value = 3259
$ws.Range("F31").Value = 371
$ws.Range("F32").Value = 752
$ws.Range("F33").Value = 4460
$ws.Range("F35").Value = 135
$ws.Range("F36").Value = 153
$ws.Range("F37").Value = 1156
$ws.Range("F38").Value = 104
$ws.Range("F41").Value = 926
$ws.Range("F42").Value = 1143
$ws.Range("F43").Value = 2061

# ---- 本地生活 sheet ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1155

# ---- 全部类型 sheet ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1155
$ws.Range("F7").Value = 3762
$ws.Range("F8").Value = 3762
$ws.Range("F9").Value = 295
$ws.Range("F10").Value = 5301
$ws.Range("F11").Value = 590
$ws.Range("F12").Value = 426
$ws.Range("F13").Value = 229
$ws.Range("F14").Value = 1061
$ws.Range("F16").Value = 141
$ws.Range("F17").Value = 50
$ws.Range("F18").Value = 731
$ws.Range("F19").Value = 364
$ws.Range("F20").Value = 47
$ws.Range("F23").Value = 179
$ws.Range("F26").Value = 6067
$ws.Range("F28").Value = 44
$ws.Range("F30").Value = 6826
$ws.Range("F33").Value = 3259
$ws.Range("F34").Value = 371
$ws.Range("F35").Value = 752
$ws.Range("F36").Value = 4460
$ws.Range("F39").Value = 135
$ws.Range("F40").Value = 153
$ws.Range("F41").Value = 1156
$ws.Range("F42").Value = 104
$ws.Range("F45").Value = 926
$ws.Range("F46").Value = 1143
$ws.Range("F48").Value = 2061
